# httk-benchmarks.xlsx: add the "2.6.1" release row to the Sheet1 / Table1
# benchmark table (R CMD check bugfix release logged after 2.6.0).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data row (row 29) -------------------------------------------------
# Match the existing table's look: left-aligned cells (style used by every
# other data row) across the whole table width, even where a column has no
# value for this release yet.
$newRow = $ws.Range("A29:R29")
$newRow.HorizontalAlignment = -4131   # xlLeft - same as the rest of the table

$ws.Range("A29").Value = "2.6.1"
$ws.Range("B29").Value = 1007
$ws.Range("C29").Value = 0.9998
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 1.001
# F29:R29 intentionally left blank - no RMSLE/notes data yet for this release.

# --- Grow the structured table so the new row is part of Table1 -----------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:R29"))

# --- Update the view state to match where the author was working ----------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C24").Select() | Out-Null
